# The "ancestor" / "tổ tiên" entry (row 4: English="ancestor", Tieng Viet="tổ tiên")
# was removed from the family vocabulary list. Deleting the entire row shifts
# the rows below it (grandfather, brother, sister) up by one and shrinks the
# used range from 7 rows to 6 rows, matching the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4").Delete()
